$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$ws = $wb.Worksheets.Item("P_valores")

$ws.Range("C2").Value = 0.9017293968466276
$ws.Range("D2").Value = 0.3138063991458253
$ws.Range("E2").Value = 0.7056299456915707
$ws.Range("F2").Value = 0.3092206450784714

$ws.Range("B3").Value = 0.9017293968466276
$ws.Range("D3").Value = 0.07654108912581759
$ws.Range("E3").Value = 0.810136114609143
$ws.Range("F3").Value = 0.496877653498959

$ws.Range("B4").Value = 0.3138063991458253
$ws.Range("C4").Value = 0.07654108912581759
$ws.Range("E4").Value = 0.09395812762150735
$ws.Range("F4").Value = 0.1626162229897858

$ws.Range("B5").Value = 0.7056299456915707
$ws.Range("C5").Value = 0.810136114609143
$ws.Range("D5").Value = 0.09395812762150735
$ws.Range("F5").Value = 0.509385455339475

$ws.Range("B6").Value = 0.3092206450784714
$ws.Range("C6").Value = 0.496877653498959
$ws.Range("D6").Value = 0.1626162229897858
$ws.Range("E6").Value = 0.509385455339475

# --- Sheet: Estadisticos_DM ---
$ws2 = $wb.Worksheets.Item("Estadisticos_DM")

$ws2.Range("C2").Value = -0.1257350965378983
$ws2.Range("D2").Value = 1.044832317283449
$ws2.Range("E2").Value = -0.385542922639723
$ws2.Range("F2").Value = -1.055162878543964

$ws2.Range("B3").Value = 0.1257350965378983
$ws2.Range("D3").Value = 1.912147129982556
$ws2.Range("E3").Value = -0.2448353336857338
$ws2.Range("F3").Value = -0.6975547831465242

$ws2.Range("B4").Value = -1.044832317283449
$ws2.Range("C4").Value = -1.912147129982556
$ws2.Range("E4").Value = -1.796859591041085
$ws2.Range("F4").Value = -1.473976744711684

$ws2.Range("B5").Value = 0.385542922639723
$ws2.Range("C5").Value = 0.2448353336857338
$ws2.Range("D5").Value = 1.796859591041085
$ws2.Range("F5").Value = -0.6770861348680258

$ws2.Range("B6").Value = 1.055162878543964
$ws2.Range("C6").Value = 0.6975547831465242
$ws2.Range("D6").Value = 1.473976744711684
$ws2.Range("E6").Value = 0.6770861348680258

$wb.Save()
